# Add an "ID" field as the new first column of Sheet1.
#
# Sheet1 currently holds: amount | name | date | reason | status  (A:E)
# with a trailing "sum" row at row 13. We insert a new column before A,
# shifting the existing data right by one column, then populate the new
# column with an "ID" header and the single data row's id value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Shift columns A:E right by inserting a fresh column at A.
[void]$ws.Columns("A:A").Insert(-4161)   # -4161 = xlShiftToRight

# New "ID" column content.
$ws.Range("A1").Value = "ID"
$ws.Range("A2").Value = 1

# The existing columns (now B:E) keep their widths, nudged to the values
# Excel recomputed for them after the insert.
$ws.Columns("B").ColumnWidth = 7.416666666666667
$ws.Columns("C").ColumnWidth = 9.541666666666666
$ws.Columns("D").ColumnWidth = 11.916666666666666
$ws.Columns("E").ColumnWidth = 14.791666666666666

# Leave the cursor where the author left it.
$ws.Range("A4").Select()
